# Property-category fix: the "building" (建物) sheet and the "car" (汽車)
# sheet rows had been tagged with property_category = "land" by mistake.
# Correct them to "building" and "car" respectively.

$wb = $excel.ActiveWorkbook

# Sheet 2 = 建物 (building). Column I is property_category.
# Rows 2 and 3 both need to say "building" instead of "land".
$wsBuilding = $wb.Worksheets.Item(2)
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"

# Sheet 3 = 汽車 (car). Column H is property_category.
# Row 2 needs to say "car" instead of "land".
$wsCar = $wb.Worksheets.Item(3)
$wsCar.Range("H2").Value = "car"
